$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.453.21"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").Value = "1.625.13"
$ws.Range("D5").Value = "212.70"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("D6").Value = "0.501"
$ws.Range("E6").Value = "  +1.36%  "
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("D10").Value = "18.71"
$ws.Range("E10").Value = "  -1.46%  "
$ws.Range("D11").Value = "0.0842"
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("D12").Value = "1.851.36"
$ws.Range("E12").Value = "  -0.60%  "
$ws.Range("D13").Value = "1.634.83"
$ws.Range("E13").Value = "  -0.02%  "
$ws.Range("D14").Value = "4.12"
$ws.Range("E14").Value = "  +1.35%  "
$ws.Range("E15").Value = "  -0.68%  "
$ws.Range("D16").Value = "64.82"
$ws.Range("E16").Value = "  +2.94%  "
$ws.Range("D17").Value = "26.492.92"
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("D19").Value = "213.64"
$ws.Range("E19").Value = "  +2.52%  "
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("E21").Value = "  -0.60%  "
$ws.Range("E22").Value = "  +1.03%  "
$ws.Range("E23").Value = "  -1.16%  "
$ws.Range("D24").Value = "2.04"
$ws.Range("E24").Value = "  +6.36%  "
$ws.Range("D25").Value = "148.45"
$ws.Range("E25").Value = "  +1.28%  "
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("E27").Value = "  -0.82%  "
$ws.Range("D28").Value = "6.84"
$ws.Range("E28").Value = "  +1.45%  "
$ws.Range("E29").Value = "  +0.70%  "
$ws.Range("D30").Value = "0.0507"
$ws.Range("E30").Value = "  -1.90%  "
$ws.Range("E31").Value = "  -0.93%  "
$ws.Range("D32").Value = "3.32"
$ws.Range("E32").Value = "  +2.79%  "
$ws.Range("E33").Value = "  -0.62%  "
$ws.Range("B34").Value = "Maker"
$ws.Range("C34").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D34").Value = "1.232.76"
$ws.Range("E34").Value = "  +5.61%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "1.49"
$ws.Range("E35").Value = "  -0.32%  "
$ws.Range("E36").Value = "  -1.88%  "
$ws.Range("D37").Value = "0.0173"
$ws.Range("E37").Value = "  +3.60%  "
$ws.Range("E38").Value = "  +0.10%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "0.506"
$ws.Range("E39").Value = "  +0.58%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "0.792"
$ws.Range("E40").Value = "  -1.89%  "
$ws.Range("E41").Value = "  -2.03%  "
$ws.Range("E42").Value = "  -0.31%  "
$ws.Range("E43").Value = "  -0.78%  "
$ws.Range("D44").Value = "1.760.50"
$ws.Range("E44").Value = "  -0.68%  "
$ws.Range("D45").Value = "92.85"
$ws.Range("E45").Value = "  +0.61%  "
$ws.Range("D46").Value = "1.58"
$ws.Range("E46").Value = "  +1.69%  "
$ws.Range("D47").Value = "54.71"
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0103"
$ws.Range("E48").Value = "  -0.65%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.0509"
$ws.Range("E49").Value = "  -0.52%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "0.406"
$ws.Range("E50").Value = "  -0.94%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "7.47"
$ws.Range("E51").Value = "  -0.62%  "
